# Applies the "Add files via upload" revision to bestellijst.xlsx:
#   - adds a new sensor row (BH1750FVI) with a name, a hyperlink and a price
#   - moves the active selection from C9 to D5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor line (row 4, next to "licht"): part number, hyperlinked part
# number, and price - mirrors the other sensor rows (6, 9, 10).
$ws.Range("B4").Value = "BH1750FVI"
$ws.Range("C4").Value = "BH1750FVI"
$ws.Range("D4").Value = "0,99$"

[void]$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.dfrobot.com/product-1713.html")
$ws.Range("C4").Style = "Hyperlink"

# Selection moved to D5 in the saved file.
[void]$ws.Range("D5").Select()
